$d = $word.ActiveDocument

for ($n = 1; $n -le 4; $n++) {
    $old = "<id>p007r_a$n</id>"
    $new = "<id>p007r_$n</id>"

    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
    }
}
